$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solicitud gráfica")
$ws.Range("C7").Value = "MA_07_05_REC10"
$ws.Range("C7").Select()
